$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge cells A1:D1 (clears the content of all but the top-left cell,
# matching openpyxl's merge_cells behaviour), then immediately unmerge
# them again. Unmerging does not restore the cleared values.
$ws.Range("A1:D1").Merge()
$ws.Range("A1:D1").UnMerge()
